$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("M11").Value = ""
$ws.Range("H19").Value = 1319.6
$ws.Range("I19").Value = 1423.7693
$ws.Range("J19").Value = 642.5
$ws.Range("K19").Value = 1423.7693
$ws.Range("L19").Value = 642.5
$ws.Range("M19").Value = -1248.7693
$ws.Range("N19").Value = -992.5
$ws.Range("H32").Value = 4736
$ws.Range("I32").Value = 5000
$ws.Range("J32").Value = 4648
$ws.Range("K32").Value = 5000
$ws.Range("L32").Value = 4648
$ws.Range("M32").Value = -4674
$ws.Range("N32").Value = -5300
$ws.Range("H38").Value = 42.125
$ws.Range("I38").Value = 42.125
$ws.Range("K38").Value = 126.375
$ws.Range("M38").Value = 245.625
$ws.Range("H43").Value = 3756.4
$ws.Range("J43").Value = 3995.5
$ws.Range("L43").Value = 3995.5
$ws.Range("N43").Value = -4133.5
$ws.Range("H51").Value = 8949.6
$ws.Range("J51").Value = 8642.857
$ws.Range("L51").Value = 8642.857
$ws.Range("N51").Value = -9610.857
$ws.Range("H103").Value = 4761.75
$ws.Range("J103").Value = 8998
$ws.Range("L103").Value = 26994
$ws.Range("N103").Value = -28166
$ws.Range("H129").Value = 20000880
$ws.Range("I129").Value = 20000880
$ws.Range("K129").Value = 60002640
$ws.Range("M129").Value = -59997640
$ws.Range("H137").Value = 2843.1428
$ws.Range("I137").Value = 2843.1428
$ws.Range("K137").Value = 8529.428400000001
$ws.Range("M137").Value = -5979.428400000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H21").Value = 7800
$ws.Range("I21").Value = 15000
$ws.Range("J21").Value = 600
$ws.Range("K21").Value = 15000
$ws.Range("L21").Value = 600
$ws.Range("M21").Value = -14626
$ws.Range("N21").Value = -1348
$ws.Range("H37").Value = 0
$ws.Range("I37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("M37").Value = ""
$ws.Range("H98").Value = 60748.25
$ws.Range("J98").Value = 60748.25
$ws.Range("L98").Value = 60748.25
$ws.Range("N98").Value = -66738.25
$ws.Range("H132").Value = 800
$ws.Range("I132").Value = 800
$ws.Range("K132").Value = 2400
$ws.Range("M132").Value = 130
$ws.Range("H139").Value = 86666.664
$ws.Range("J139").Value = 86666.664
$ws.Range("L139").Value = 86666.664
$ws.Range("N139").Value = -96946.664

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 363
$ws.Range("I22").Value = 293.5
$ws.Range("K22").Value = 293.5
$ws.Range("M22").Value = -120.5
$ws.Range("H132").Value = 124000
$ws.Range("I132").Value = 124000
$ws.Range("K132").Value = 124000
$ws.Range("M132").Value = -118940
$ws.Range("H134").Value = 3000
$ws.Range("I134").Value = 2500
$ws.Range("K134").Value = 7500
$ws.Range("M134").Value = -4965
$ws.Range("H140").Value = 98999
$ws.Range("J140").Value = 98999
$ws.Range("L140").Value = 98999
$ws.Range("N140").Value = -109359

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 209.44444
$ws.Range("I22").Value = 198
$ws.Range("J22").Value = 223.75
$ws.Range("K22").Value = 198
$ws.Range("L22").Value = 223.75
$ws.Range("M22").Value = 152
$ws.Range("N22").Value = -923.75
$ws.Range("H28").Value = 34500
$ws.Range("J28").Value = 34500
$ws.Range("L28").Value = 34500
$ws.Range("N28").Value = -34990
$ws.Range("H31").Value = 11119.875
$ws.Range("I31").Value = 7247.25
$ws.Range("J31").Value = 14992.5
$ws.Range("K31").Value = 7247.25
$ws.Range("L31").Value = 14992.5
$ws.Range("M31").Value = -6952.25
$ws.Range("N31").Value = -15582.5
$ws.Range("H32").Value = 3050
$ws.Range("I32").Value = 1350
$ws.Range("J32").Value = 4750
$ws.Range("K32").Value = 1350
$ws.Range("L32").Value = 4750
$ws.Range("M32").Value = -1034
$ws.Range("N32").Value = -5382
$ws.Range("H34").Value = 11119.875
$ws.Range("I34").Value = 7247.25
$ws.Range("J34").Value = 14992.5
$ws.Range("K34").Value = 7247.25
$ws.Range("L34").Value = 14992.5
$ws.Range("M34").Value = -7045.25
$ws.Range("N34").Value = -15396.5
$ws.Range("H102").Value = 64998
$ws.Range("J102").Value = 64998
$ws.Range("L102").Value = 64998
$ws.Range("N102").Value = -69866
$ws.Range("H105").Value = 2030.3334
$ws.Range("I105").Value = 2030.3334
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 2030.3334
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = ""
$ws.Range("N105").Value = -283.3334
$ws.Range("H134").Value = 13506
$ws.Range("I134").Value = 10013
$ws.Range("J134").Value = 16999
$ws.Range("K134").Value = 30039
$ws.Range("L134").Value = 50997
$ws.Range("M134").Value = -27504
$ws.Range("N134").Value = -56067

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 31.461538
$ws.Range("I12").Value = 56.5
$ws.Range("J12").Value = 20.333334
$ws.Range("K12").Value = 169.5
$ws.Range("L12").Value = 61.000002
$ws.Range("M12").Value = 3.5
$ws.Range("N12").Value = -407.000002
$ws.Range("H80").Value = 1000
$ws.Range("I80").Value = 1000
$ws.Range("K80").Value = 3000
$ws.Range("M80").Value = -2064
$ws.Range("H83").Value = 1000
$ws.Range("I83").Value = 1000
$ws.Range("K83").Value = 9000
$ws.Range("M83").Value = -4320
$ws.Range("H92").Value = 14998
$ws.Range("J92").Value = 14998
$ws.Range("L92").Value = 44994
$ws.Range("N92").Value = -47490

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H47").Value = 35000
$ws.Range("J47").Value = 35000
$ws.Range("L47").Value = 35000
$ws.Range("N47").Value = -36136
$ws.Range("H80").Value = 7627
$ws.Range("I80").Value = 7627
$ws.Range("K80").Value = 7627
$ws.Range("M80").Value = -6629
$ws.Range("H83").Value = 7627
$ws.Range("I83").Value = 7627
$ws.Range("K83").Value = 38135
$ws.Range("M83").Value = -33143

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("L2").Value = ""
$ws.Range("N2").Value = 0
$ws.Range("H32").Value = 1310
$ws.Range("I32").Value = 1310
$ws.Range("K32").Value = 1310
$ws.Range("M32").Value = -993
$ws.Range("H132").Value = 9996.23
$ws.Range("I132").Value = 9996
$ws.Range("K132").Value = 29988
$ws.Range("M132").Value = -27458

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H39").Value = 7000
$ws.Range("I39").Value = 7000
$ws.Range("K39").Value = 7000
$ws.Range("M39").Value = -6587
$ws.Range("H94").Value = 45666.668
$ws.Range("J94").Value = 45666.668
$ws.Range("L94").Value = 45666.668
$ws.Range("N94").Value = -47468.668
$ws.Range("H100").Value = 1626
$ws.Range("I100").Value = 1179.2
$ws.Range("K100").Value = 2358.4
$ws.Range("M100").Value = -1817.4
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = ""
$ws.Range("N103").Value = 0
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = ""
$ws.Range("N105").Value = 0
$ws.Range("H113").Value = 2480.111
$ws.Range("J113").Value = 2970.5
$ws.Range("L113").Value = 8911.5
$ws.Range("N113").Value = -13251.5
$ws.Range("H136").Value = 3399.6667
$ws.Range("I136").Value = 3399.6667
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 10199.0001
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = ""
$ws.Range("N136").Value = -7649.000100000001
